$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-10: full cross-join of Sending cluster (A) x Target cluster (D)
# over {ECs, FAPs, sCs}; B (Nlgn3) / C (Nrxn2) constant throughout.

# row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Nlgn3"
$ws.Cells.Item(2,3).Value = "Nrxn2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.3239616666666667
$ws.Cells.Item(2,8).Value = 0.9718850000000001
$ws.Cells.Item(2,9).Value = 0.1555281019885789
$ws.Cells.Item(2,10).Value = 0.1555281019885789
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.039512
$ws.Cells.Item(2,14).Value = 0.118536
$ws.Cells.Item(2,15).Value = 0.006183359004302676
$ws.Cells.Item(2,16).Value = 0.006183359004302677
$ws.Cells.Item(2,17).Value = 0.01280037337333334
$ws.Cells.Item(2,18).Value = 0.11520336036
$ws.Cells.Item(2,19).Value = 0.0009616860898531841
$ws.Cells.Item(2,20).Value = 0.0009616860898531841

# row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Nlgn3"
$ws.Cells.Item(3,3).Value = "Nrxn2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.3239616666666667
$ws.Cells.Item(3,8).Value = 0.9718850000000001
$ws.Cells.Item(3,9).Value = 0.1555281019885789
$ws.Cells.Item(3,10).Value = 0.1555281019885789
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 6.039054
$ws.Cells.Item(3,14).Value = 18.117162
$ws.Cells.Item(3,15).Value = 0.9450708374258476
$ws.Cells.Item(3,16).Value = 0.9450708374258477
$ws.Cells.Item(3,17).Value = 1.95642199893
$ws.Cells.Item(3,18).Value = 17.60779799037
$ws.Cells.Item(3,19).Value = 0.1469850735895989
$ws.Cells.Item(3,20).Value = 0.1469850735895989

# row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Nlgn3"
$ws.Cells.Item(4,3).Value = "Nrxn2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.3239616666666667
$ws.Cells.Item(4,8).Value = 0.9718850000000001
$ws.Cells.Item(4,9).Value = 0.1555281019885789
$ws.Cells.Item(4,10).Value = 0.1555281019885789
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.3114883333333333
$ws.Cells.Item(4,14).Value = 0.934465
$ws.Cells.Item(4,15).Value = 0.04874580356984966
$ws.Cells.Item(4,16).Value = 0.04874580356984967
$ws.Cells.Item(4,17).Value = 0.1009102796138889
$ws.Cells.Item(4,18).Value = 0.9081925165250001
$ws.Cells.Item(4,19).Value = 0.007581342309126811
$ws.Cells.Item(4,20).Value = 0.007581342309126811

# row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Nlgn3"
$ws.Cells.Item(5,3).Value = "Nrxn2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.026996666666667
$ws.Cells.Item(5,8).Value = 3.08099
$ws.Cells.Item(5,9).Value = 0.4930424144274184
$ws.Cells.Item(5,10).Value = 0.4930424144274184
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.039512
$ws.Cells.Item(5,14).Value = 0.118536
$ws.Cells.Item(5,15).Value = 0.006183359004302676
$ws.Cells.Item(5,16).Value = 0.006183359004302677
$ws.Cells.Item(5,17).Value = 0.04057869229333333
$ws.Cells.Item(5,18).Value = 0.36520823064
$ws.Cells.Item(5,19).Value = 0.003048658252752909
$ws.Cells.Item(5,20).Value = 0.003048658252752909

# row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Nlgn3"
$ws.Cells.Item(6,3).Value = "Nrxn2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.026996666666667
$ws.Cells.Item(6,8).Value = 3.08099
$ws.Cells.Item(6,9).Value = 0.4930424144274184
$ws.Cells.Item(6,10).Value = 0.4930424144274184
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 6.039054
$ws.Cells.Item(6,14).Value = 18.117162
$ws.Cells.Item(6,15).Value = 0.9450708374258476
$ws.Cells.Item(6,16).Value = 0.9450708374258477
$ws.Cells.Item(6,17).Value = 6.202088327819999
$ws.Cells.Item(6,18).Value = 55.81879495038
$ws.Cells.Item(6,19).Value = 0.4659600074893822
$ws.Cells.Item(6,20).Value = 0.4659600074893822

# row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Nlgn3"
$ws.Cells.Item(7,3).Value = "Nrxn2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.026996666666667
$ws.Cells.Item(7,8).Value = 3.08099
$ws.Cells.Item(7,9).Value = 0.4930424144274184
$ws.Cells.Item(7,10).Value = 0.4930424144274184
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.3114883333333333
$ws.Cells.Item(7,14).Value = 0.934465
$ws.Cells.Item(7,15).Value = 0.04874580356984966
$ws.Cells.Item(7,16).Value = 0.04874580356984967
$ws.Cells.Item(7,17).Value = 0.3198974800388888
$ws.Cells.Item(7,18).Value = 2.87907732035
$ws.Cells.Item(7,19).Value = 0.02403374868528335
$ws.Cells.Item(7,20).Value = 0.02403374868528335

# row 8: sCs -> ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Nlgn3"
$ws.Cells.Item(8,3).Value = "Nrxn2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.73202
$ws.Cells.Item(8,8).Value = 2.19606
$ws.Cells.Item(8,9).Value = 0.3514294835840027
$ws.Cells.Item(8,10).Value = 0.3514294835840027
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.039512
$ws.Cells.Item(8,14).Value = 0.118536
$ws.Cells.Item(8,15).Value = 0.006183359004302676
$ws.Cells.Item(8,16).Value = 0.006183359004302677
$ws.Cells.Item(8,17).Value = 0.02892357424
$ws.Cells.Item(8,18).Value = 0.26031216816
$ws.Cells.Item(8,19).Value = 0.002173014661696583
$ws.Cells.Item(8,20).Value = 0.002173014661696583

# row 9: sCs -> FAPs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Nlgn3"
$ws.Cells.Item(9,3).Value = "Nrxn2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.73202
$ws.Cells.Item(9,8).Value = 2.19606
$ws.Cells.Item(9,9).Value = 0.3514294835840027
$ws.Cells.Item(9,10).Value = 0.3514294835840027
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 6.039054
$ws.Cells.Item(9,14).Value = 18.117162
$ws.Cells.Item(9,15).Value = 0.9450708374258476
$ws.Cells.Item(9,16).Value = 0.9450708374258477
$ws.Cells.Item(9,17).Value = 4.42070830908
$ws.Cells.Item(9,18).Value = 39.78637478172001
$ws.Cells.Item(9,19).Value = 0.3321257563468666
$ws.Cells.Item(9,20).Value = 0.3321257563468666

# row 10: sCs -> sCs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Nlgn3"
$ws.Cells.Item(10,3).Value = "Nrxn2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.73202
$ws.Cells.Item(10,8).Value = 2.19606
$ws.Cells.Item(10,9).Value = 0.3514294835840027
$ws.Cells.Item(10,10).Value = 0.3514294835840027
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.3114883333333333
$ws.Cells.Item(10,14).Value = 0.934465
$ws.Cells.Item(10,15).Value = 0.04874580356984966
$ws.Cells.Item(10,16).Value = 0.04874580356984967
$ws.Cells.Item(10,17).Value = 0.2280156897666666
$ws.Cells.Item(10,18).Value = 2.0521412079
$ws.Cells.Item(10,19).Value = 0.0171307125754395
$ws.Cells.Item(10,20).Value = 0.01713071257543951

